$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''57.397.29'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '''3.102.93'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''524.59'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').Value = '''136.55'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.39%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '''3.104.13'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +2.27%  '
$ws.Range('D10').Value = '''7.23'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.38%  '
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').Value = '''3.644.23'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('E14').Value = '  +2.88%  '
$ws.Range('E15').Value = '  -3.37%  '
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '''57.491.80'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '''3.103.74'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '''5.92'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.82%  '
$ws.Range('D20').Value = '''12.43'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.73%  '
$ws.Range('D21').Value = '''7.86'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('D22').Value = '''347.04'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.56%  '
$ws.Range('E23').Value = '  -0.83%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').Value = '''68.18'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +2.17%  '
$ws.Range('E26').Value = '  -2.23%  '
$ws.Range('E27').Value = '  -0.86%  '
$ws.Range('D28').Value = '''1.00'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').Value = '''0.0₃0904'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.94%  '
$ws.Range('D30').Value = '''0.998'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = '''7.37'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.54%  '
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').Value = '''5.98'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -8.10%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('D36').Value = '''4.90'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +6.01%  '
$ws.Range('D37').Value = '''157.94'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.33%  '
$ws.Range('D38').Value = '''6.09'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.30%  '
$ws.Range('D39').Value = '''25.92'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.41%  '
$ws.Range('E40').Value = '  -3.54%  '
$ws.Range('D41').Value = '''4.17'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +5.83%  '
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('E43').Value = '  +5.79%  '
$ws.Range('D44').Value = '''0.698'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.59%  '
$ws.Range('D45').Value = '''3.138.57'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Value = '''36.42'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('D47').Value = '''0.999'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').Value = '''2.338.15'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.38%  '
$ws.Range('E49').Value = '  +2.90%  '
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').Value = '''0.948'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.55%  '
